# edit.ps1 - apply the CV content update described by the commit diff
# Uses Word COM-interop Find + Range.Text against $word.ActiveDocument.
#
# We locate each old string with Range.Find.Execute(FindText) (search-only,
# no ReplaceWith) and then assign the new string straight to Range.Text.
# Doing the substitution this way (rather than passing ReplaceWith to
# Find.Execute) avoids the AutoCorrect/AutoFormat "smart quotes" pass that
# Word's real Replace pipeline applies, so straight apostrophes in the new
# text (e.g. "Webpack's") are preserved verbatim.

$d = $word.ActiveDocument

function Replace-Text {
    param(
        [string]$OldText,
        [string]$NewText
    )
    $rng = $d.Content
    $ok = $rng.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false)
    if ($ok) {
        $rng.Text = $NewText
    } else {
        Write-Output "NOT FOUND: $OldText"
    }
    return $ok
}

# 1. Title under the name: "SOFTWARE ENGINEER" -> "Software Engineer"
Replace-Text "SOFTWARE ENGINEER" "Software Engineer"

# 2. Summary paragraph rewrite
Replace-Text `
    "Experienced Frontend Engineer expertise in JavaScript, ReactJS, and Web Technologies, complemented by a strong background in API Design, API management. An Open-Source enthusiast. Passionate about creating exceptional user experiences." `
    "Experienced Frontend Software Engineer with expertise in JavaScript, ReactJS and Web Technologies, complemented by a strong background in API Design, API management. An Open-Source enthusiast. Passionate about creating exceptional user experiences."

# 3. Programming Languages: drop "Java"
Replace-Text `
    "JavaScript | Typescript | Python | HTML5 | CSS3 | Java | PHP | Ruby" `
    "JavaScript | Typescript | Python | HTML5 | CSS3 | PHP | Ruby"

# 4. Front-end skills list update
Replace-Text `
    " | Cypress | Styled-components | Redux | Accessibility | WCAG | React-i18n | Material-UI" `
    " | Styled-components | Redux | Microfrontends | WCAG | Figma | Material-UI"

# 5. Back-end skills: add Kubernetes
Replace-Text `
    "ExpressJS | NestJS | Python Flask | Python Django | Firebase | AWS EC2" `
    "ExpressJS | NestJS | Python Flask | Python Django | Firebase | AWS EC2 | Kubernetes"

# 6. Application Security: "OpenID Connect" -> "OpenID"
Replace-Text `
    "OpenID Connect | JWT | CORS | OAuth2 | OWASP | Content Security Policy" `
    "OpenID| JWT | CORS | OAuth2 | OWASP | Content Security Policy"

# 7. AWS role header: merge location text and update the employment dates
Replace-Text ", Toronto, Canada" ", Toronto, Canada"
Replace-Text "March 2023-Current" "March 2023-June 2023"

# 8. AWS bullet: End-to-End testing detail expanded
Replace-Text `
    "Added a set of End-toEnd tests using AWS CloudWatch Synthetics." `
    "Added a set of End-to-End tests using AWS CloudWatch Synthetics for the Simple Notification Service to verify the functionality in the Edit subscription flow in the SNS console."

# 9. OAuth PKCE bullet: trim the trailing metric
Replace-Text `
    "Developed a SPA authentication architecture by implementing the OAuth PKCE extension flow, resulting in significantly enhanced security and a performance improvement of about 50%." `
    "Developed a SPA authentication architecture by implementing the OAuth PKCE extension flow, resulting in significantly enhanced security and a performance improvement "

# 10. State management bullet rewrite
Replace-Text `
    "Designed and developed a simple state management mechanism using context and reducer hooks for Choreo SPA application, resulting in a significant improvement in user experience." `
    "simplified the state management architecture in the Choreo SPA console by using context and reducer hooks, thereby replacing some of the Redux usage. "

# 11. Design system bullet: normalize to a single run (text unchanged)
Replace-Text `
    "Developed a design system based on Material-UI that significantly enhanced the developer experience and ensured consistent product UI across the entire application." `
    "Developed a design system based on Material-UI that significantly enhanced the developer experience and ensured consistent product UI across the entire application."

# 12. Associate Technical Lead bullet: loading/performance rewrite
Replace-Text `
    "Improved SPA performance by incorporating optimization techniques such as lazy loading, strong compression, memorization." `
    "Improved the initial loading time and performance of the Choreo portal by utilizing Webpack's dynamic imports and React's lazy loading techniques to chunk the JS bundles."

# 13. React SDK bullet -> caching bullet rewrite
Replace-Text `
    "Developed a React (Typescript) SDK to communicate with product REST APIs, enabling efficient data retrieval and manipulation." `
    "Improved the static files and JS bundles caching by incorporating Webpack's chunk hashes and content hashing mechanisms, which eliminates the occurrence of stale content for users."

Write-Output "Done"
